$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "34.688.60"
Set-TextValue $ws.Range("E2") "  +0.57%  "
Set-TextValue $ws.Range("D3") "1.840.80"
Set-TextValue $ws.Range("E3") "  +1.75%  "
Set-TextValue $ws.Range("E4") "  +0.31%  "
Set-TextValue $ws.Range("D5") "227.19"
Set-TextValue $ws.Range("E5") "  +0.60%  "
Set-TextValue $ws.Range("E6") "  +2.03%  "
Set-TextValue $ws.Range("E7") "  +0.22%  "
Set-TextValue $ws.Range("D8") "42.49"
Set-TextValue $ws.Range("E8") "  +17.27%  "
Set-TextValue $ws.Range("D9") "0.303"
Set-TextValue $ws.Range("E9") "  +3.59%  "
Set-TextValue $ws.Range("D10") "0.0686"
Set-TextValue $ws.Range("E10") "  +0.87%  "
Set-TextValue $ws.Range("E11") "  +3.52%  "
Set-TextValue $ws.Range("D12") "2.109.54"
Set-TextValue $ws.Range("E12") "  +1.85%  "
Set-TextValue $ws.Range("D13") "11.30"
Set-TextValue $ws.Range("E13") "  +0.20%  "
Set-TextValue $ws.Range("D14") "1.823.87"
Set-TextValue $ws.Range("E14") "  +0.28%  "
Set-TextValue $ws.Range("D15") "4.70"
Set-TextValue $ws.Range("E15") "  +6.40%  "
Set-TextValue $ws.Range("E16") "  +3.54%  "
Set-TextValue $ws.Range("D17") "34.702.36"
Set-TextValue $ws.Range("E17") "  +0.68%  "
Set-TextValue $ws.Range("D18") "68.62"
Set-TextValue $ws.Range("E18") "  +0.32%  "
Set-TextValue $ws.Range("D19") "243.42"
Set-TextValue $ws.Range("E19") "  +0.47%  "
Set-TextValue $ws.Range("D20") "0.0₃0785"
Set-TextValue $ws.Range("E20") "  +1.32%  "
Set-TextValue $ws.Range("D21") "12.10"
Set-TextValue $ws.Range("E21") "  +7.95%  "
Set-TextValue $ws.Range("D22") "4.72"
Set-TextValue $ws.Range("E22") "  +15.37%  "
Set-TextValue $ws.Range("E23") "  +0.23%  "
Set-TextValue $ws.Range("E24") "  -1.62%  "
Set-TextValue $ws.Range("D25") "172.11"
Set-TextValue $ws.Range("E25") "  +0.37%  "
Set-TextValue $ws.Range("D26") "7.97"
Set-TextValue $ws.Range("E26") "  +1.18%  "
Set-TextValue $ws.Range("D27") "17.79"
Set-TextValue $ws.Range("E27") "  +3.27%  "
Set-TextValue $ws.Range("E28") "  +0.53%  "
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  +0.34%  "
Set-TextValue $ws.Range("E30") "  +5.36%  "
Set-TextValue $ws.Range("E31") "  +2.35%  "
Set-TextValue $ws.Range("D32") "4.00"
Set-TextValue $ws.Range("E32") "  +2.31%  "
Set-TextValue $ws.Range("D33") "0.0528"
Set-TextValue $ws.Range("E33") "  +2.04%  "
Set-TextValue $ws.Range("D34") "1.87"
Set-TextValue $ws.Range("E34") "  +4.37%  "
Set-TextValue $ws.Range("D35") "90.02"
Set-TextValue $ws.Range("E35") "  +11.32%  "
Set-TextValue $ws.Range("D36") "0.663"
Set-TextValue $ws.Range("E36") "  +1.60%  "
Set-TextValue $ws.Range("D37") "1.342.33"
Set-TextValue $ws.Range("E37") "  -1.54%  "
Set-TextValue $ws.Range("D38") "2.44"
Set-TextValue $ws.Range("E38") "  +3.82%  "
Set-TextValue $ws.Range("E39") "  +0.72%  "
Set-TextValue $ws.Range("E40") "  +3.66%  "
Set-TextValue $ws.Range("D41") "14.97"
Set-TextValue $ws.Range("E41") "  +11.89%  "
Set-TextValue $ws.Range("D42") "0.985"
Set-TextValue $ws.Range("E42") "  +5.20%  "
Set-TextValue $ws.Range("E43") "  +6.87%  "
Set-TextValue $ws.Range("D44") "2.82"
Set-TextValue $ws.Range("E44") "  +1.47%  "
Set-TextValue $ws.Range("D45") "2.44"
Set-TextValue $ws.Range("E45") "  +0.59%  "
Set-TextValue $ws.Range("E46") "  +4.16%  "
Set-TextValue $ws.Range("D47") "2.008.12"
Set-TextValue $ws.Range("E47") "  +1.86%  "
Set-TextValue $ws.Range("E48") "  +3.26%  "
Set-TextValue $ws.Range("E49") "  +0.20%  "
Set-TextValue $ws.Range("D50") "102.39"
Set-TextValue $ws.Range("E50") "  -0.02%  "
Set-TextValue $ws.Range("D51") "0.0613"
Set-TextValue $ws.Range("E51") "  +1.10%  "
